$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: round values to 2 decimal places (matching the target dataset)
$ws.Range("B5").Value = 7.69
$ws.Range("C5").Value = 5.44
$ws.Range("D5").Value = 0.8100000000000001
$ws.Range("E5").Value = 16.42
$ws.Range("F5").Value = 13.6
$ws.Range("G5").Value = 6.05
$ws.Range("H5").Value = 21.56
$ws.Range("I5").Value = 9.31
$ws.Range("J5").Value = 3.98
$ws.Range("K5").Value = 6.07
$ws.Range("L5").Value = 6.64
$ws.Range("M5").Value = 6.87
$ws.Range("N5").Value = 1.93
$ws.Range("O5").Value = 6.02
$ws.Range("P5").Value = 8.41
$ws.Range("Q5").Value = 5.21
$ws.Range("R5").Value = 0.76
$ws.Range("S5").Value = 0.45
$ws.Range("T5").Value = 84.18000000000001
$ws.Range("U5").Value = 16.63
$ws.Range("V5").Value = 5.55
$ws.Range("W5").Value = 10.97
$ws.Range("X5").Value = 6.01
$ws.Range("Y5").Value = 0.71
$ws.Range("Z5").Value = 10.17
$ws.Range("AA5").Value = 4.9
$ws.Range("AB5").Value = 4.45
$ws.Range("AC5").Value = 5.21
$ws.Range("AD5").Value = 6.93
$ws.Range("AE5").Value = 0.5600000000000001
$ws.Range("AF5").Value = 18.97
$ws.Range("AG5").Value = 3.06
$ws.Range("AH5").Value = 6.94

# Delete row 6 entirely (data trimmed)
$ws.Rows.Item(6).Delete()
